$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 6-8: third "Groupe" (3) + three more students.
$ws.Range("A6").Value = 3
$ws.Range("B6").Value = "Truc Much"
$ws.Range("B7").Value = "Car Resse"
$ws.Range("B8").Value = "Aloe Vera"

# Row-height tweaks present in the target sheet.
$ws.Range("A3").RowHeight = 13.8
$ws.Range("A6").RowHeight = 12.8

# Move the active selection to C15.
$ws.Range("C15").Select()

# Workbook window tab-ratio (Calc -> 72%). Best effort: engine may not
# persist this particular window-chrome property.
$excel.ActiveWindow.TabRatio = 0.72
